# Updated cryptos list with GitHub Actions: refresh Price (D) and Volume(1h) (E)
# columns for the coin table on the active worksheet.
#
# The Price column stores values as plain text (e.g. "60.050.12",
# "0.0\u20830768") which are not valid Excel numbers, so each Price cell is
# assigned with a leading single-quote to force a text value and avoid
# Excel re-interpreting/mangling the string as a number or date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '60.050.12'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = "'" + '2.396.92'
$ws.Range("E3").Value = '  -0.96%  '
$ws.Range("D5").Value = "'" + '559.03'
$ws.Range("E5").Value = '  +1.19%  '
$ws.Range("D6").Value = "'" + '134.61'
$ws.Range("E6").Value = '  -2.12%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = "'" + '0.587'
$ws.Range("E8").Value = '  -0.21%  '
$ws.Range("D9").Value = "'" + '0.105'
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("D10").Value = "'" + '5.63'
$ws.Range("E10").Value = '  -1.03%  '
$ws.Range("E11").Value = '  +1.37%  '
$ws.Range("E12").Value = '  -2.60%  '
$ws.Range("D13").Value = "'" + '24.59'
$ws.Range("E13").Value = '  -3.68%  '
$ws.Range("D14").Value = "'" + '2.824.97'
$ws.Range("E14").Value = '  -0.86%  '
$ws.Range("D15").Value = "'" + '59.964.89'
$ws.Range("E15").Value = '  +0.14%  '
$ws.Range("E16").Value = '  -0.43%  '
$ws.Range("D17").Value = "'" + '2.396.59'
$ws.Range("E17").Value = '  -1.39%  '
$ws.Range("D18").Value = "'" + '11.13'
$ws.Range("E18").Value = '  -2.18%  '
$ws.Range("D19").Value = "'" + '4.51'
$ws.Range("E19").Value = '  +2.23%  '
$ws.Range("D20").Value = "'" + '322.58'
$ws.Range("E20").Value = '  -2.17%  '
$ws.Range("D21").Value = "'" + '6.75'
$ws.Range("E21").Value = '  +1.44%  '
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").Value = "'" + '64.28'
$ws.Range("E23").Value = '  -3.35%  '
$ws.Range("D24").Value = "'" + '0.174'
$ws.Range("E24").Value = '  +0.98%  '
$ws.Range("D25").Value = "'" + '8.50'
$ws.Range("E25").Value = '  -1.77%  '
$ws.Range("D26").Value = "'" + '1.00'
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("E27").Value = '  +0.57%  '
$ws.Range("E28").Value = '  +1.91%  '
$ws.Range("D29").Value = "'" + '0.0₃0768'
$ws.Range("E29").Value = '  -1.30%  '
$ws.Range("D30").Value = "'" + '171.11'
$ws.Range("E30").Value = '  +1.10%  '
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("E32").Value = '  +7.54%  '
$ws.Range("D33").Value = "'" + '0.400'
$ws.Range("E33").Value = '  -2.30%  '
$ws.Range("D34").Value = "'" + '18.26'
$ws.Range("E34").Value = '  -2.25%  '
$ws.Range("E36").Value = '  +2.16%  '
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("E38").Value = '  -1.68%  '
$ws.Range("D39").Value = "'" + '322.00'
$ws.Range("E39").Value = '  +2.50%  '
$ws.Range("E40").Value = '  -0.80%  '
$ws.Range("D41").Value = "'" + '38.70'
$ws.Range("E41").Value = '  -2.28%  '
$ws.Range("D42").Value = "'" + '146.74'
$ws.Range("E42").Value = '  +5.87%  '
$ws.Range("E43").Value = '  -3.47%  '
$ws.Range("D44").Value = "'" + '0.0969'
$ws.Range("E44").Value = '  +0.22%  '
$ws.Range("D45").Value = "'" + '19.81'
$ws.Range("E45").Value = '  +1.16%  '
$ws.Range("E46").Value = '  -1.29%  '
$ws.Range("E47").Value = '  -0.91%  '
$ws.Range("E48").Value = '  -1.90%  '
$ws.Range("D49").Value = "'" + '11.06'
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("E51").Value = '  +0.10%  '
